$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly records were added to the daily log for
# "Terminal Hortofrutícola Agro Chillán" / Durazno / Florida King.
# They belong right after the existing row 102, so insert two blank
# rows there (pushing the old rows 103.. down to 105..) and fill them in.

$ws.Rows.Item(103).Insert()
$ws.Rows.Item(103).Insert()

# --- New row 103: Florida King / Primera ---
$ws.Cells.Item(103, 1).Value = 7
$ws.Cells.Item(103, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(103, 3).Value = "Ñuble"
$ws.Cells.Item(103, 4).Value = 44529
$ws.Cells.Item(103, 5).Value = 16
$ws.Cells.Item(103, 6).Value = "Fruta"
$ws.Cells.Item(103, 7).Value = 100103
$ws.Cells.Item(103, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(103, 9).Value = 100103004
$ws.Cells.Item(103, 10).Value = "Durazno"
$ws.Cells.Item(103, 11).Value = "Florida King"
$ws.Cells.Item(103, 12).Value = "Primera"
$ws.Cells.Item(103, 13).Value = 120
$ws.Cells.Item(103, 14).Value = 18000
$ws.Cells.Item(103, 15).Value = 19000
$ws.Cells.Item(103, 16).Value = 18500
$ws.Cells.Item(103, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(103, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(103, 19).Value = 1156
$ws.Cells.Item(103, 20).Value = 16

# --- New row 104: Florida King / Segunda ---
$ws.Cells.Item(104, 1).Value = 7
$ws.Cells.Item(104, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(104, 3).Value = "Ñuble"
$ws.Cells.Item(104, 4).Value = 44529
$ws.Cells.Item(104, 5).Value = 16
$ws.Cells.Item(104, 6).Value = "Fruta"
$ws.Cells.Item(104, 7).Value = 100103
$ws.Cells.Item(104, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(104, 9).Value = 100103004
$ws.Cells.Item(104, 10).Value = "Durazno"
$ws.Cells.Item(104, 11).Value = "Florida King"
$ws.Cells.Item(104, 12).Value = "Segunda"
$ws.Cells.Item(104, 13).Value = 80
$ws.Cells.Item(104, 14).Value = 16000
$ws.Cells.Item(104, 15).Value = 16000
$ws.Cells.Item(104, 16).Value = 16000
$ws.Cells.Item(104, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(104, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(104, 19).Value = 1000
$ws.Cells.Item(104, 20).Value = 16
